$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new 19:00 timeslot entry for 2026/02/23 (月) was added to the log.
# This pushes row 858 ("2026/12/29") and every row after it down by one,
# growing the sheet from A1:D899 to A1:D900.
$ws.Rows(858).Insert(-4121)

# Fill in the newly inserted row 858 with the new data point.
# Force the date column to text first so the "yyyy/mm/dd" string is stored
# as a literal string (matching the rest of the column) instead of being
# auto-parsed into a date serial number; then restore the plain style so no
# stray number-format style is left behind on the cell.
$ws.Cells.Item(858, 1).NumberFormat = "@"
$ws.Cells.Item(858, 1).Value = "2026/02/23"
$ws.Cells.Item(858, 1).Style = $ws.Cells.Item(857, 1).Style

$ws.Cells.Item(858, 2).Value = "月"
$ws.Cells.Item(858, 3).Value = 19
$ws.Cells.Item(858, 4).Value = 201
